$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2:H2").NumberFormat = "@"

$ws.Range("C2").Value = "{'criterion': 'gini', 'max_depth': 5, 'min_samples_leaf': 4, 'min_samples_split': 2}"
$ws.Range("E2").Value = "58.99%"
$ws.Range("F2").Value = "60.27%"
$ws.Range("G2").Value = "58.99%"
$ws.Range("H2").Value = "58.42%"
